$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '51.564.42'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '2.987.75'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '382.37'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").Value = '103.61'
$ws.Range("E6").Value = '  +2.04%  '
$ws.Range("D7").Value = '0.546'
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").Value = '36.74'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '0.0859'
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = '3.442.29'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '18.44'
$ws.Range("E14").Value = '  +1.84%  '
$ws.Range("D15").Value = '7.81'
$ws.Range("E15").Value = '  +2.63%  '
$ws.Range("D16").Value = '2.982.70'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '11.18'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '0.999'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '51.502.30'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '3.08'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '12.63'
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").Value = '70.45'
$ws.Range("E23").Value = '  +2.45%  '
$ws.Range("D24").Value = '267.54'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = '3.23'
$ws.Range("E25").Value = '  +2.05%  '
$ws.Range("D26").Value = '7.91'
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("D27").Value = '7.33'
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("E28").Value = '  +2.83%  '
$ws.Range("D30").Value = '26.12'
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").Value = '10.34'
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("D33").Value = '34.79'
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").Value = '51.59'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").Value = '0.0442'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  +3.52%  '
$ws.Range("D39").Value = '16.80'
$ws.Range("E39").Value = '  +2.38%  '
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").Value = '2.56'
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("D43").Value = '126.32'
$ws.Range("E43").Value = '  +4.72%  '
$ws.Range("E44").Value = '  +11.44%  '
$ws.Range("D45").Value = '21.50'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").Value = '2.37'
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("D48").Value = '0.271'
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("D49").Value = '2.030.21'
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("D50").Value = '0.0334'
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("E51").Value = '  +15.56%  '
